$wb = $excel.ActiveWorkbook

# Sheet "展览" (first sheet) - update column F ("想去人数") values
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 1001
$wsExhibit.Range("F4").Value = 161
$wsExhibit.Range("F5").Value = 2705
$wsExhibit.Range("F7").Value = 208
$wsExhibit.Range("F9").Value = 109
$wsExhibit.Range("F11").Value = 2530
$wsExhibit.Range("F12").Value = 637

# Sheet "全部类型" (fourth sheet) - update column F ("想去人数") values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1002
$wsAll.Range("F5").Value = 161
$wsAll.Range("F6").Value = 2705
$wsAll.Range("F8").Value = 208
$wsAll.Range("F11").Value = 109
$wsAll.Range("F13").Value = 2530
$wsAll.Range("F14").Value = 637
